$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 54) below the existing data (which ends at row 53).
# The date column stores plain text like "2025/10/03" (not a real date value),
# so we prefix with an apostrophe to force text entry and avoid Excel's
# automatic date-serial conversion; then reset the style to "Normal" so no
# stray number-format style is left attached to the cell.
$ws.Cells.Item(54, 1).Value = "'2025/10/03"
$ws.Cells.Item(54, 1).Style = "Normal"

$ws.Cells.Item(54, 2).Value = "金"
$ws.Cells.Item(54, 3).Value = 9
$ws.Cells.Item(54, 4).Value = 3
